# Increase MCGLT to hit 1.5 target
# Bump the max capacity growth value from 500 to 700 across the MCGLT
# lookup table (column C), and update the one formula-driven cell (C26)
# to reference the new value. Column D cells are formulas referencing
# column C, so they recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCGLT")

$rows = @(2, 4, 6, 8, 10, 12, 14, 16, 18, 20, 22, 24, 28, 30, 32)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 700
}

# C26 holds a formula that mirrors the 500 -> 700 bump.
$ws.Range("C26").Formula = "=IF(About!C2=1,700,0)"

# Match the author's recorded active-cell selection on the MCGLT sheet,
# then restore "About" as the active tab (it was the active tab before
# this edit and stays that way afterwards).
$ws.Range("C33").Select()
$wb.Worksheets.Item("About").Activate()
